$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text formatting so that
# values such as "580.70" or "67.568.78" are not re-interpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.568.78"
$ws.Range("E2").Value = "  -0.03%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.324.14"
$ws.Range("E3").Value = "  +0.31%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.11%  "

# Row 5 - BNB
$ws.Range("D5").Value = "580.70"
$ws.Range("E5").Value = "  -0.58%  "

# Row 6 - Solana
$ws.Range("D6").Value = "175.34"
$ws.Range("E6").Value = "  -3.66%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.05%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  +0.00%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "3.321.95"
$ws.Range("E9").Value = "  +0.40%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.21%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +0.10%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "45.21"
$ws.Range("E12").Value = "  -2.02%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -1.45%  "

# Row 14 - BitcoinCash
$ws.Range("D14").Value = "666.31"
$ws.Range("E14").Value = "  +4.49%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.868.38"
$ws.Range("E15").Value = "  +0.52%  "

# Row 16 - Polkadot
$ws.Range("D16").Value = "8.39"
$ws.Range("E16").Value = "  -0.36%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "67.588.98"
$ws.Range("E17").Value = "  -0.23%  "

# Row 18 - now TRON (was WrappedEther)
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.118"
$ws.Range("E18").Value = "  -0.72%  "

# Row 19 - now WrappedEther (was TRON)
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.325.20"
$ws.Range("E19").Value = "  +0.16%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  -1.17%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "10.95"
$ws.Range("E21").Value = "  +0.68%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "0.886"
$ws.Range("E22").Value = "  -1.45%  "

# Row 23 - Toncoin
$ws.Range("D23").Value = "5.38"
$ws.Range("E23").Value = "  +7.23%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "17.01"
$ws.Range("E24").Value = "  -3.06%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "98.64"
$ws.Range("E25").Value = "  +1.93%  "

# Row 26 - PancakeSwap
$ws.Range("E26").Value = "  -4.03%  "

# Row 27 - ImmutableX
$ws.Range("E27").Value = "  -3.63%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "9.24"
$ws.Range("E28").Value = "  -3.35%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "33.46"
$ws.Range("E29").Value = "  +2.76%  "

# Row 30 - Filecoin
$ws.Range("E30").Value = "  -1.41%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "7.29"
$ws.Range("E31").Value = "  +9.70%  "

# Row 32 - Bittensor
$ws.Range("D32").Value = "569.60"
$ws.Range("E32").Value = "  -3.20%  "

# Row 33 - Cosmos
$ws.Range("D33").Value = "10.93"
$ws.Range("E33").Value = "  +0.20%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  +0.78%  "

# Row 35 - Dai
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.01%  "

# Row 36 - now OKB (was Maker)
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "56.63"
$ws.Range("E36").Value = "  +1.99%  "

# Row 37 - now Maker (was OKB)
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "3.675.00"
$ws.Range("E37").Value = "  -6.56%  "

# Row 38 - dogwifhat
$ws.Range("D38").Value = "3.29"
$ws.Range("E38").Value = "  -6.31%  "

# Row 39 - InjectiveProtocol
$ws.Range("D39").Value = "34.08"
$ws.Range("E39").Value = "  +5.00%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  +0.29%  "

# Row 41 - Fetch.AI
$ws.Range("E41").Value = "  -1.81%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  -4.24%  "

# Row 43 - ApeXProtocol
$ws.Range("D43").Value = "3.33"
$ws.Range("E43").Value = "  -1.58%  "

# Row 44 - TheGraph
$ws.Range("E44").Value = "  -0.99%  "

# Row 45 - PEPE
$ws.Range("E45").Value = "  -2.74%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  -1.94%  "

# Row 47 - ThetaToken
$ws.Range("E47").Value = "  +2.24%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  -0.09%  "

# Row 49 - FirstDigitalUSD
$ws.Range("E49").Value = "  -0.42%  "

# Row 50 - Mantle
$ws.Range("D50").Value = "1.36"
$ws.Range("E50").Value = "  -1.06%  "

# Row 51 - Monero
$ws.Range("D51").Value = "129.32"
$ws.Range("E51").Value = "  -1.06%  "
